$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 542, shifting existing rows 542:612 down to 543:613
$ws.Rows(542).Insert()

# Populate the new row 542 with its data
$ws.Range("A542").Value = 10
$ws.Range("B542").Value = "Vega Modelo de Temuco"
$ws.Range("C542").Value = "La Araucanía"
$ws.Range("D542").Value = 45127
$ws.Range("E542").Value = 9
$ws.Range("F542").Value = "Fruta"
$ws.Range("G542").Value = 100108
$ws.Range("H542").Value = "Tropicales y subtropicales"
$ws.Range("I542").Value = 100108002
$ws.Range("J542").Value = "Mango"
$ws.Range("K542").Value = "Sin especificar"
$ws.Range("L542").Value = "Primera"
$ws.Range("M542").Value = 780
$ws.Range("N542").Value = 9000
$ws.Range("O542").Value = 9000
$ws.Range("P542").Value = 9000
$ws.Range("Q542").Value = "$/bandeja 4 kilos"
$ws.Range("R542").Value = "Brasil"
$ws.Range("S542").Value = 2250
$ws.Range("T542").Value = 4
